$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add value 456 to A2
$ws.Range("A2").Value = 456

# Move the active selection to A3 (reflects the cursor having moved past A2)
$ws.Range("A3").Select()
